$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is a brand new data row for "Crumpet" (GEF / exporter).
# Row 6 already exists in the template as a blank, pre-styled row; we only
# need to populate its values (its styles already match rows 2-4).

# --- Row 5: Crumpet GEF ---
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 5 is new, so it does not yet carry the shared formatting that rows
# 2-4 (and the blank row 6) already have. Copy that formatting across
# without touching values that were just written.
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 6: Scone GEF ---
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"
